# "Man with a plan" -- replace the old axes/panel placeholder labels on
# Sheet3 with the new "Final plan" outline, and restore the original
# selections/active sheet afterwards.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws3 = $wb.Worksheets.Item(3)

# --- Sheet3: clear the old placeholder data (axes1/panel8/.../panel21 and
# the stray numbers) so the shared strings they used can be recycled ---
$ws3.Range("A3").ClearContents()
$ws3.Range("B3").ClearContents()
$ws3.Range("C3").ClearContents()
$ws3.Range("D3").ClearContents()
$ws3.Range("E3").ClearContents()
$ws3.Range("B4").ClearContents()
$ws3.Range("C4").ClearContents()
$ws3.Range("D4").ClearContents()
$ws3.Range("C5").ClearContents()

# --- Sheet3: type in the new "Final plan" outline, in entry order so the
# shared-string table comes out in the same order Excel would produce it ---
$ws3.Range("A2").Value = "Final plan:"
$ws3.Range("A3").Value = "crop"
$ws3.Range("B3").Value = "gray"
$ws3.Range("C3").Value = "hist"
$ws3.Range("D3").Value = "OC50"
$ws3.Range("E3").Value = "edge (what kind?)"
$ws3.Range("F3").Value = "D10E15"
$ws3.Range("F4").Value = "can tune"
$ws3.Range("G3").Value = "leave only largest body"
$ws3.Range("H3").Value = "fill with random points"
$ws3.Range("I3").Value = "find circle "
$ws3.Range("J3").Value = "estimate circle size"

# --- Sheet3 picked up an explicit (portrait) page setup, matching the
# other two sheets in the workbook ---
$ws3.PageSetup.Orientation = 1

# --- restore the selections recorded in the saved workbook: Sheet1's
# cursor moved to L7, Sheet3's to G4, and Sheet3 is the tab left showing ---
[void]$ws1.Range("L7").Select()
[void]$ws3.Select()
[void]$ws3.Range("G4").Select()
